$d = $word.ActiveDocument
$xml = '<w:p><w:r><w:t xml:space="preserve">1. </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>訓練給定中</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Microsoft JhengHei UI" w:eastAsia="Microsoft JhengHei UI" w:hAnsi="Microsoft JhengHei UI" w:cs="Microsoft JhengHei UI" w:hint="eastAsia"/></w:rPr><w:t>⽂</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="新細明體" w:eastAsia="新細明體" w:hAnsi="新細明體" w:cs="新細明體" w:hint="eastAsia"/></w:rPr><w:t>文語料庫</w:t></w:r><w:r><w:t xml:space="preserve"> (1.1G) 520</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>萬筆資料</w:t></w:r><w:r><w:t xml:space="preserve"> (JSON)</w:t></w:r></w:p><w:p><w:r><w:t>https://drive.google.com/open?</w:t></w:r></w:p><w:p><w:r><w:t>id=1EX8eE5YWBxCaohBO8Fh4e2j3b9C2bTVQ</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">2. </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>利利</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Microsoft JhengHei UI" w:eastAsia="Microsoft JhengHei UI" w:hAnsi="Microsoft JhengHei UI" w:cs="Microsoft JhengHei UI" w:hint="eastAsia"/></w:rPr><w:t>⽤</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="新細明體" w:eastAsia="新細明體" w:hAnsi="新細明體" w:cs="新細明體" w:hint="eastAsia"/></w:rPr><w:t>用</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>keras</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> / </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pytorch</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> / seq2seq-lstm / </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>sklearn</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">seq2seq </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>英翻中實作</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">3. </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>輸入</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Microsoft JhengHei UI" w:eastAsia="Microsoft JhengHei UI" w:hAnsi="Microsoft JhengHei UI" w:cs="Microsoft JhengHei UI" w:hint="eastAsia"/></w:rPr><w:t>⼀</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="新細明體" w:eastAsia="新細明體" w:hAnsi="新細明體" w:cs="新細明體" w:hint="eastAsia"/></w:rPr><w:t>一個簡單的英</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Microsoft JhengHei UI" w:eastAsia="Microsoft JhengHei UI" w:hAnsi="Microsoft JhengHei UI" w:cs="Microsoft JhengHei UI" w:hint="eastAsia"/></w:rPr><w:t>⽂</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="新細明體" w:eastAsia="新細明體" w:hAnsi="新細明體" w:cs="新細明體" w:hint="eastAsia"/></w:rPr><w:t>文句句</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Microsoft JhengHei UI" w:eastAsia="Microsoft JhengHei UI" w:hAnsi="Microsoft JhengHei UI" w:cs="Microsoft JhengHei UI" w:hint="eastAsia"/></w:rPr><w:t>⼦</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="新細明體" w:eastAsia="新細明體" w:hAnsi="新細明體" w:cs="新細明體" w:hint="eastAsia"/></w:rPr><w:t>子</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">3.1 </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>“</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>It is a nice day</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>”</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>—</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">&gt; </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>今天天氣真好</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">3.2 </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>最後</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>100</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>筆為</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>testing data (</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>需呈現</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>testing</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>結果</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>)</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">3.3 </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>輸出中</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Microsoft JhengHei UI" w:eastAsia="Microsoft JhengHei UI" w:hAnsi="Microsoft JhengHei UI" w:cs="Microsoft JhengHei UI" w:hint="eastAsia"/></w:rPr><w:t>⽂</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="新細明體" w:eastAsia="新細明體" w:hAnsi="新細明體" w:cs="新細明體" w:hint="eastAsia"/></w:rPr><w:t>文翻譯結果</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">3.4 </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>分數以輸出中</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Microsoft JhengHei UI" w:eastAsia="Microsoft JhengHei UI" w:hAnsi="Microsoft JhengHei UI" w:cs="Microsoft JhengHei UI" w:hint="eastAsia"/></w:rPr><w:t>⽂</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="新細明體" w:eastAsia="新細明體" w:hAnsi="新細明體" w:cs="新細明體" w:hint="eastAsia"/></w:rPr><w:t>文翻譯結果評</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>分</w:t></w:r></w:p>'
$target = $d.Paragraphs.Last

if ($target.Range.Text.Trim().Length -eq 0) {
    $target.Range.InsertXML($xml)
} else {
    $target.Range.InsertParagraphAfter()
    $target = $d.Paragraphs.Last
    $target.Range.InsertXML($xml)
}
